# Add columns I ("I0") and J ("IF") to the sheet, mirroring the style of
# the existing header cells and filling in the per-row data:
#   I = 1 for every data row
#   J = copy of column H for every data row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from H1 onto I1:J1,
# then set the header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the data rows 2-26.
for ($r = 2; $r -le 26; $r++) {
    $hValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hValue
}
